# "Check all links performance tuning,File upload"
# Adds a new "File upload" row (B23/C23) to the object repository sheet,
# mirroring the existing "find by id" / "fileToUpload" row (row 19),
# and extends the "select one" dropdown validation down to the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 19 (Find By method = "find by id", Value = "fileToUpload")
# into the new row 23, picking up matching cell formatting.
$ws.Range("B19:C19").Copy($ws.Range("B23:C23"))
$ws.Range("B23").Value = "find by id"
$ws.Range("C23").Value = "fileToUpload"

# Extend the "select one" list validation that covered B21:B22 so it also
# covers the newly added B23 cell.
$ws.Range("B21:B22").Validation.Delete()
$validation = $ws.Range("B21:B23").Validation
$validation.Add(3, 1, 1, """find by id,find by xpath,find by name,find by css,find by linktext,find by partial linktext,find by tagname""")
$validation.InputTitle = "select one"
$validation.IgnoreBlank = $true
$validation.InCellDropdown = $true
$validation.ShowInput = $true
$validation.ShowError = $true

# Match the author's final selection/navigation state.
$ws.Range("B23").Select()
